# Auto-generated Excel COM script to apply the Ravana Profits sheet update
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific rows
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H9").Value = 927.36365
$ws.Range("I9").Value = 1061.3334
$ws.Range("J9").Value = 324.5
$ws.Range("K9").Value = 1061.3334
$ws.Range("L9").Value = 324.5
$ws.Range("M9").Value = -892.3334
$ws.Range("N9").Value = -662.5
$ws.Range("H12").Value = 50
$ws.Range("J12").Value = 50
$ws.Range("L12").Value = 50
$ws.Range("N12").Value = -390
$ws.Range("H19").Value = 1849.9
$ws.Range("I19").Value = 1542.4286
$ws.Range("J19").Value = 2567.3333
$ws.Range("K19").Value = 1542.4286
$ws.Range("L19").Value = 2567.3333
$ws.Range("M19").Value = -1367.4286
$ws.Range("N19").Value = -2917.3333
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H28").Value = 949.75
$ws.Range("I28").Value = 999.5
$ws.Range("K28").Value = 999.5
$ws.Range("M28").Value = -514.5
$ws.Range("H29").Value = 1900
$ws.Range("I29").Value = 50
$ws.Range("K29").Value = 150
$ws.Range("M29").Value = 131
$ws.Range("H32").Value = 8929.75
$ws.Range("J32").Value = 8573
$ws.Range("L32").Value = 8573
$ws.Range("N32").Value = -9225
$ws.Range("H38").Value = 1667.75
$ws.Range("I38").Value = 139
$ws.Range("K38").Value = 417
$ws.Range("M38").Value = -45
$ws.Range("H121").Value = 1837.4286
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1837.4286
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 5512.2858
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -9006.2858
$ws.Range("H137").Value = 2975.111
$ws.Range("I137").Value = 1111.1538
$ws.Range("J137").Value = 4705.9287
$ws.Range("K137").Value = 3333.4614
$ws.Range("L137").Value = 14117.7861
$ws.Range("M137").Value = -783.4614000000001
$ws.Range("N137").Value = -19217.7861

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H61").Value = 6981.3335
$ws.Range("I61").Value = 7778.2
$ws.Range("K61").Value = 7778.2
$ws.Range("M61").Value = -7566.2
$ws.Range("H74").Value = 1215.909
$ws.Range("I74").Value = 1322.25
$ws.Range("J74").Value = 932.3333
$ws.Range("K74").Value = 1322.25
$ws.Range("L74").Value = 932.3333
$ws.Range("M74").Value = -448.25
$ws.Range("N74").Value = -2680.3333
$ws.Range("H77").Value = 1215.909
$ws.Range("I77").Value = 1322.25
$ws.Range("J77").Value = 932.3333
$ws.Range("K77").Value = 6611.25
$ws.Range("L77").Value = 4661.6665
$ws.Range("M77").Value = -2243.25
$ws.Range("N77").Value = -13397.6665
$ws.Range("H102").Value = 3462.2
$ws.Range("I102").Value = 2674
$ws.Range("J102").Value = 4644.5
$ws.Range("K102").Value = 2674
$ws.Range("L102").Value = 4644.5
$ws.Range("M102").Value = -1052
$ws.Range("N102").Value = -7888.5
$ws.Range("H132").Value = 1547.8096
$ws.Range("I132").Value = 1290.8108
$ws.Range("K132").Value = 3872.4324
$ws.Range("M132").Value = -1342.4324
$ws.Range("H136").Value = 6981.3335
$ws.Range("I136").Value = 7778.2
$ws.Range("K136").Value = 23334.6
$ws.Range("M136").Value = -20784.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4619.5
$ws.Range("J86").Value = 4833.3335
$ws.Range("L86").Value = 4833.3335
$ws.Range("N86").Value = -7079.3335
$ws.Range("H89").Value = 4619.5
$ws.Range("J89").Value = 4833.3335
$ws.Range("L89").Value = 24166.6675
$ws.Range("N89").Value = -35398.6675
$ws.Range("H105").Value = 3530.1667
$ws.Range("I105").Value = 3530.1667
$ws.Range("K105").Value = 3530.1667
$ws.Range("M105").Value = -1783.1667
$ws.Range("H107").Value = 1157
$ws.Range("I107").Value = 1157
$ws.Range("K107").Value = 1157
$ws.Range("M107").Value = 763
$ws.Range("H134").Value = 3467.4443
$ws.Range("I134").Value = 3149.5625
$ws.Range("J134").Value = 6010.5
$ws.Range("K134").Value = 9448.6875
$ws.Range("L134").Value = 18031.5
$ws.Range("M134").Value = -6913.6875
$ws.Range("N134").Value = -23101.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 629
$ws.Range("I22").Value = 309.66666
$ws.Range("J22").Value = 820.6
$ws.Range("K22").Value = 309.66666
$ws.Range("L22").Value = 820.6
$ws.Range("M22").Value = 40.33334000000002
$ws.Range("N22").Value = -1520.6
$ws.Range("H58").Value = 2273.625
$ws.Range("I58").Value = 2184.2856
$ws.Range("J58").Value = 2343.111
$ws.Range("K58").Value = 2184.2856
$ws.Range("L58").Value = 2343.111
$ws.Range("M58").Value = -1981.2856
$ws.Range("N58").Value = -2749.111
$ws.Range("I105").Value = 5000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3253
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2124.4285
$ws.Range("I122").Value = 2149.75
$ws.Range("J122").Value = 1972.5
$ws.Range("K122").Value = 6449.25
$ws.Range("L122").Value = 5917.5
$ws.Range("M122").Value = -3999.25
$ws.Range("N122").Value = -10817.5
$ws.Range("H132").Value = 3101.8635
$ws.Range("I132").Value = 2828.5789
$ws.Range("J132").Value = 4832.6665
$ws.Range("K132").Value = 8485.736699999999
$ws.Range("L132").Value = 14497.9995
$ws.Range("M132").Value = -5955.736699999999
$ws.Range("N132").Value = -19557.9995
$ws.Range("H136").Value = 2273.625
$ws.Range("I136").Value = 2184.2856
$ws.Range("J136").Value = 2343.111
$ws.Range("K136").Value = 6552.8568
$ws.Range("L136").Value = 7029.333
$ws.Range("M136").Value = -4002.8568
$ws.Range("N136").Value = -12129.333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 81.2
$ws.Range("J2").Value = 81.2
$ws.Range("L2").Value = 487.2
$ws.Range("N2").Value = -713.2
$ws.Range("H38").Value = 281.2
$ws.Range("I38").Value = 251
$ws.Range("J38").Value = 402
$ws.Range("K38").Value = 753
$ws.Range("L38").Value = 1206
$ws.Range("M38").Value = -406
$ws.Range("N38").Value = -1900
$ws.Range("H107").Value = 263
$ws.Range("J107").Value = 280.77777
$ws.Range("L107").Value = 842.33331
$ws.Range("N107").Value = -4682.33331

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 76.333336
$ws.Range("I2").Value = 91.666664
$ws.Range("J2").Value = 45.666668
$ws.Range("K2").Value = 91.666664
$ws.Range("L2").Value = 45.666668
$ws.Range("M2").Value = 21.333336
$ws.Range("N2").Value = -271.666668
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 3513.1428
$ws.Range("I132").Value = 2918.8
$ws.Range("K132").Value = 8756.400000000001
$ws.Range("M132").Value = -6226.400000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H106").Value = 17500
$ws.Range("J106").Value = 17500
$ws.Range("L106").Value = 17500
$ws.Range("N106").Value = -20024
$ws.Range("H136").Value = 3105.6875
$ws.Range("I136").Value = 3046.0667
$ws.Range("K136").Value = 9138.2001
$ws.Range("M136").Value = -6588.2001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11323.333
$ws.Range("I81").Value = 6484
$ws.Range("K81").Value = 12968
$ws.Range("M81").Value = -11907
$ws.Range("H84").Value = 11323.333
$ws.Range("I84").Value = 6484
$ws.Range("K84").Value = 64840
$ws.Range("M84").Value = -59536
$ws.Range("H132").Value = 1974.8823
$ws.Range("I132").Value = 1315.6957
$ws.Range("K132").Value = 3947.0871
$ws.Range("M132").Value = -1417.0871
$ws.Range("H136").Value = 825
$ws.Range("I136").Value = 646.9231
$ws.Range("K136").Value = 1940.7693
$ws.Range("M136").Value = 609.2307000000001
